$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Team" and "Claims_Officer_Name" rows (rows 4 and 5)
$ws.Range("A4:C5").EntireRow.Delete()

# Remove the trailing "Entitlement_Weeks" row (now row 19 after the above delete)
$ws.Range("A19:C19").EntireRow.Delete()

# Leave the selection on A8, matching the saved view state
$ws.Range("A8").Select()
